# Update return logic and fix code flow
# Adds new registration-form columns (H:L) to the header row and appends
# a sample registration record in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (H1:L1) -------------------------------------------
$ws.Range("H1").Value = "Blood Donation"
$ws.Range("I1").Value = "Blood Group"
$ws.Range("J1").Value = "Webinar Interest"
$ws.Range("K1").Value = "Webinar Date"
$ws.Range("L1").Value = "Registered At"

# --- New data row (A2:L2) -------------------------------------------------
$ws.Range("A2").Value = "Good Tester"

# WhatsApp number must stay plain text (not be coerced to a number), so mark
# the cell as text before assigning the value.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "9111111111"

$ws.Range("C2").Value = "good@test.local"
$ws.Range("D2").Value = "PG"
$ws.Range("E2").Value = "IT Professional"
$ws.Range("F2").Value = "Female"
$ws.Range("G2").Value = "Good College"
$ws.Range("H2").Value = "Yes"
$ws.Range("I2").Value = "A+"
$ws.Range("J2").Value = "Yes"

# Webinar date / registration timestamp are stored as plain text strings in
# the source data, not as real Excel dates, so force text storage here too.
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2025-12-31"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2025-12-24 16:58:59"
